# Commit message: "renamed repo, fixed output folder path"
# The underlying data change is that duplicate rows (caused by the old
# output path bug) were removed from two of the per-species result sheets.
#
# Sheet "s__Fenollaria sp900539225-b-p" (1st sheet): remove rows 9-11
#   (label_UMGS110_11.fasta, label_UMGS110_18.fasta, label_UMGS110_2.fasta),
#   shrinking the used range from A1:F32 to A1:F29.
#
# Sheet "s__Fenollaria sp900539725-b-p" (2nd sheet): remove rows 8-18
#   (label_UMGS167_12.fasta .. label_UMGS167_6.fasta),
#   shrinking the used range from A1:F51 to A1:F40.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows("9:11").Delete()

$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows("8:18").Delete()
